$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking price strings (e.g. '65.255.81')
# that must remain plain text, matching the source data exactly. Force the
# cell to Text format before assigning, then restore the default style so
# no stray formatting is introduced.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '65.255.81'
$ws.Range('E2').Value = '  -1.13%  '
Set-TextValue 'D3' '2.948.26'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('E4').Value = '  -0.14%  '
Set-TextValue 'D5' '569.71'
$ws.Range('E5').Value = '  -2.23%  '
Set-TextValue 'D6' '160.93'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  -0.04%  '
Set-TextValue 'D8' '0.516'
$ws.Range('E8').Value = '  -0.30%  '
Set-TextValue 'D9' '2.944.82'
$ws.Range('E9').Value = '  -2.21%  '
Set-TextValue 'D10' '6.68'
$ws.Range('E10').Value = '  -4.08%  '
$ws.Range('E11').Value = '  -3.79%  '
Set-TextValue 'D12' '0.459'
$ws.Range('E12').Value = '  +0.71%  '
Set-TextValue 'D13' '0.0000244'
$ws.Range('E13').Value = '  -2.71%  '
Set-TextValue 'D14' '34.54'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('E15').Value = '  -1.44%  '
Set-TextValue 'D16' '65.276.08'
$ws.Range('E16').Value = '  -1.20%  '
Set-TextValue 'D17' '3.410.96'
$ws.Range('E17').Value = '  -2.94%  '
Set-TextValue 'D18' '7.04'
$ws.Range('E18').Value = '  +1.32%  '
Set-TextValue 'D19' '2.944.04'
$ws.Range('E19').Value = '  -2.26%  '
Set-TextValue 'D20' '15.77'
$ws.Range('E20').Value = '  +13.25%  '
Set-TextValue 'D21' '444.14'
$ws.Range('E21').Value = '  -3.02%  '
Set-TextValue 'D22' '0.697'
$ws.Range('E22').Value = '  +1.87%  '
Set-TextValue 'D23' '7.29'
$ws.Range('E23').Value = '  -0.38%  '
Set-TextValue 'D24' '82.42'
$ws.Range('E24').Value = '  +0.20%  '
Set-TextValue 'D26' '12.18'
$ws.Range('E26').Value = '  -1.46%  '
$ws.Range('E27').Value = '  +0.03%  '
Set-TextValue 'D28' '10.01'
$ws.Range('E28').Value = '  -5.67%  '
Set-TextValue 'D29' '2.47'
$ws.Range('E29').Value = '  +5.07%  '
Set-TextValue 'D30' '7.94'
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('E31').Value = '  -0.67%  '
Set-TextValue 'D32' '0.0000102'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('E33').Value = '  +2.34%  '
Set-TextValue 'D34' '27.18'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('E35').Value = '  -0.05%  '
Set-TextValue 'D36' '0.972'
$ws.Range('E36').Value = '  -1.85%  '
Set-TextValue 'D37' '5.73'
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D38' '49.14'
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('B39').Value = 'Arweave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue 'D39' '44.96'
$ws.Range('E39').Value = '  +2.86%  '
Set-TextValue 'D40' '1.99'
$ws.Range('E40').Value = '  -7.17%  '
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D43' '2.82'
$ws.Range('E43').Value = '  -5.37%  '
$ws.Range('B44').Value = 'Cosmos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D44' '8.56'
$ws.Range('E44').Value = '  +0.74%  '
Set-TextValue 'D45' '381.37'
$ws.Range('E45').Value = '  -0.62%  '
Set-TextValue 'D46' '0.0350'
$ws.Range('E46').Value = '  -1.63%  '
Set-TextValue 'D47' '2.678.39'
$ws.Range('E47').Value = '  -3.90%  '
Set-TextValue 'D48' '133.42'
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D50' '23.65'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D51' '2.17'
$ws.Range('E51').Value = '  +1.48%  '
